# Update "Last Updated" timestamp on the Metadata sheet.
$wb = $excel.ActiveWorkbook
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("A2").Value = "05 Nov 2025, 01:09 PM"

# The "Stock List" sheet refreshed: every existing row's data moved up by
# one position (row N now shows what used to be in row N+1) and a brand
# new row was appended at the bottom with the latest stock.
$ws = $wb.Worksheets.Item("Stock List")

$firstDataRow = 2
$lastDataRow = 76

# Snapshot the current contents of rows (firstDataRow+1 .. lastDataRow)
# before overwriting anything, since the shift reads from row r+1.
$snapshot = @{}
for ($r = $firstDataRow + 1; $r -le $lastDataRow; $r++) {
    $snapshot[$r] = @{
        B = $ws.Cells.Item($r, 2).Value()
        C = $ws.Cells.Item($r, 3).Value()
        D = $ws.Cells.Item($r, 4).Value()
        E = $ws.Cells.Item($r, 5).Value()
        H = $ws.Cells.Item($r, 8).Value()
    }
}

# Shift rows up: row r becomes what row r+1 used to hold.
for ($r = $firstDataRow; $r -le $lastDataRow - 1; $r++) {
    $src = $snapshot[$r + 1]
    $ws.Cells.Item($r, 2).Value = $src.B
    $ws.Cells.Item($r, 3).Value = $src.C
    $ws.Cells.Item($r, 4).Value = $src.D
    $ws.Cells.Item($r, 5).Value = $src.E
    $ws.Cells.Item($r, 8).Value = $src.H
}

# New last row: TRAVELFOOD enters the list.
$ws.Cells.Item($lastDataRow, 2).Value = "TRAVELFOOD"
$ws.Cells.Item($lastDataRow, 3).Value = "TRAVELFOOD"
$ws.Cells.Item($lastDataRow, 4).Value = 1316.3
$ws.Cells.Item($lastDataRow, 5).Value = 0.1141
$ws.Cells.Item($lastDataRow, 8).Value = 17332.9705
